# Apply the "automatic update" edit to the Översikt KNIVSTA workbook:
#  1. Refresh the "Förändrad" (changed) date in column C for every existing
#     data row (2-143) from 2023-09-19 (45188) to 2023-09-20 (45189).
#  2. Give row 143 an explicit row height (matches every other data row).
#  3. Append a brand new felling-notification record as row 144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "changed" timestamp column for all existing rows (2-143).
$ws.Range("C2:C143").Value = 45189

# 2. Row 143 picks up an explicit custom height, like the surrounding rows.
$ws.Rows.Item(143).RowHeight = 15

# 3. Append the new row 144 with the additional felling notification.
$ws.Range("A144").Value = "A 43908-2023"

$ws.Range("B144").Value = 45187
$ws.Range("B144").NumberFormat = "YYYY-MM-DD"

$ws.Range("C144").Value = 45189
$ws.Range("C144").NumberFormat = "YYYY-MM-DD"

$ws.Range("D144").Value = "UPPSALA LÄN"
$ws.Range("E144").Value = "KNIVSTA"
$ws.Range("F144").Value = "Kyrkan"

$ws.Range("G144").Value = 1.7
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 0
$ws.Range("J144").Value = 0
$ws.Range("K144").Value = 0
$ws.Range("L144").Value = 0
$ws.Range("M144").Value = 0
$ws.Range("N144").Value = 0
$ws.Range("O144").Value = 0
$ws.Range("P144").Value = 0
$ws.Range("Q144").Value = 0

# R144 stays textually empty but keeps the wrap-text formatting used
# throughout column R.
$ws.Range("R144").WrapText = $true
